$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BOM update: filter Z1 changed from a 0402-package part (BLM15PX601SN1D) to the
# same 0603-package part already used by Z2 (BLM18HE152SN1D). The two separate
# BOM lines for Z1 and Z2 therefore collapse into a single line "Z1, Z2" with
# quantity 2, keeping Z2's Comment/Description/Part Number/Manufacturer.
#
# Row 51 currently holds Z1 (BLM15PX601SN1D, no Manufacturer) and row 52 holds
# Z2 (BLM18HE152SN1D, Murata). Deleting row 51 removes the obsolete Z1 line and
# shifts Z2 (and everything below it, e.g. ZQ1/ZQ2) up by one row, preserving
# all cell formatting along the way.
$ws.Rows("51").Delete()

# Row 51 now contains the former Z2 data. Bump the quantity to 2 and relabel
# the designator to cover both Z1 and Z2.
$ws.Range("A51").Value = 2
$ws.Range("B51").Value = "Z1, Z2"

# Re-stamp B51 with the same cell formatting (border + text style) as its
# neighbour C51, since assigning .Value resets the style to the default.
$ws.Range("C51").Copy()
$ws.Range("B51").PasteSpecial(-4122)
$excel.CutCopyMode = $false
